$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "สูตรเรขาคณิตจำยากจัง"
$ws.Range("B6").Value = "เรขาคณิต"
$ws.Range("C6").Value = "คณิตศาสตร์"

$ws.Range("A7").Value = "ลูกแมวสามารถให้อาหารอะไรได้บ้าง"
$ws.Range("B7").Value = "แมว"
$ws.Range("C7").Value = "สัตว์เลี้ยง"
$ws.Range("D7").Value = "สัตว์"

$ws.Range("A8").Value = "ทำอย่างไรแมวถึงจะยอมอาบน้ำ"
$ws.Range("B8").Value = "แมว"
$ws.Range("C8").Value = "สัตว์เลี้ยง"
$ws.Range("D8").Value = "สัตว์"
$ws.Range("E8").Value = "อาบน้ำแมว"

$ws.Range("A9").Value = "เผลอซักผ้าสีปนกับผ้าขาว ทำอย่างไรดี"
$ws.Range("B9").Value = "ซักผ้า"
$ws.Range("C9").Value = "งานบ้าน"
$ws.Range("D9").Value = "ผ้าสี"
$ws.Range("E9").Value = "ผ้าขาว"
